$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column L: header + break_on_off flag values (1 at rows 19, 37, 54; else 0)
$ws.Range("L1").Value = "break_on_off"

for ($r = 2; $r -le 73; $r++) {
    if ($r -eq 19 -or $r -eq 37 -or $r -eq 54) {
        $ws.Cells.Item($r, 12).Value = 1
    } else {
        $ws.Cells.Item($r, 12).Value = 0
    }
}

$ws.Range("L1:L73").Select()
